# 3.1.20.docx — "Added function to write attachments to txt file."
#
# 1) The table's column grid (w:tblGrid) is out of sync with the actual
#    cell widths (w:tcW) it holds; nudging each column's Width (even to
#    its current value) makes Word resynchronise w:tblGrid to the cell
#    widths, which is the change the diff shows.
# 2) After the table, add an "AND/OR" paragraph, a blank paragraph, and a
#    new paragraph that documents the file with missing correspondence
#    parties being written out to an attachment (3.1.20.txt), with the
#    placeholder token bolded/underlined.

$d = $word.ActiveDocument

# --- 1. Resync the table's grid widths with its cell widths ---
$tbl = $d.Tables(1)
$tbl.Columns(1).Width = $tbl.Columns(1).Width
$tbl.Columns(2).Width = $tbl.Columns(2).Width
$tbl.Columns(4).Width = $tbl.Columns(4).Width

# --- 2. Append the new paragraphs after the paragraph that follows the table ---
$anchor = $d.Paragraphs.Last.Range
$anchor.InsertParagraphAfter()
$pAndOr = $d.Paragraphs($d.Paragraphs.Count)
$pAndOr.Range.InsertAfter("AND/OR")

$pAndOr.Range.InsertParagraphAfter()
$pBlank = $d.Paragraphs($d.Paragraphs.Count)

$pBlank.Range.InsertParagraphAfter()
$pDetail = $d.Paragraphs($d.Paragraphs.Count)

$r = $pDetail.Range
$r.InsertAfter("Det er ikke likt antall korrespondanseparter som journalposter, da ")

# -1 because Range.End on a paragraph range sits past the trailing
# paragraph mark; the real end of the just-typed text is End - 1.
$termStart = $r.End - 1
$r.InsertAfter("ANTALLREGISTRERINGERUTENKORRESPONDANSEPART")
$termEnd = $r.End - 1
$term = $d.Range($termStart, $termEnd)
$term.Font.Bold = 1
$term.Font.Underline = 1

$r.InsertAfter(" registreringer mangler dette. ")
$r.InsertAfter("Oversikt over j")
$r.InsertAfter("ournalpostene som mangler dette")
$r.InsertAfter(" finnes i vedlegget «3.1.20.txt».")
